$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.372.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.89%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.59%  '

$ws.Range("E6").Value = '  -1.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4528'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3865'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07910'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.017'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.872.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.914'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.118'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.25%  '

$ws.Range("E17").Value = '  -3.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '85.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06521'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.39%  '

$ws.Range("E21").Value = '  -1.18%  '

$ws.Range("E22").Value = '  -4.85%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.373.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.280'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.079.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.071'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.440'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.486'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09296'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9348'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.600'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.248'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02234'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.223'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.41%  '

$ws.Range("E39").Value = '  -2.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.216'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.92%  '

$ws.Range("E41").Value = '  -1.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5900'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1887'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.279'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5615'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.367'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.921'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06768'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '108.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.71%  '
